# Update cryptos list figures (price + 1h volume change) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.439.69"
$ws.Range("E2").Value = "  +2.68%  "
$ws.Range("D3").Value = "3.193.23"
$ws.Range("E3").Value = "  +1.50%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'597.56"
$ws.Range("E5").Value = "  +3.30%  "
$ws.Range("D6").Value = "'155.01"
$ws.Range("E6").Value = "  +4.23%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +6.12%  "
$ws.Range("D9").Value = "3.193.62"
$ws.Range("E9").Value = "  +1.58%  "
$ws.Range("E10").Value = "  +1.75%  "
$ws.Range("D11").Value = "'5.91"
$ws.Range("E11").Value = "  -3.26%  "
$ws.Range("E12").Value = "  +3.76%  "
$ws.Range("E13").Value = "  +2.68%  "
$ws.Range("D14").Value = "'39.32"
$ws.Range("E14").Value = "  +6.08%  "
$ws.Range("D15").Value = "3.713.18"
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("D16").Value = "'7.51"
$ws.Range("E16").Value = "  +5.35%  "
$ws.Range("D17").Value = "66.425.66"
$ws.Range("E17").Value = "  +2.50%  "
$ws.Range("D18").Value = "3.192.21"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("D20").Value = "'519.96"
$ws.Range("E20").Value = "  +3.21%  "
$ws.Range("D21").Value = "'15.45"
$ws.Range("E21").Value = "  +3.90%  "
$ws.Range("E22").Value = "  +3.66%  "
$ws.Range("E23").Value = "  +5.47%  "
$ws.Range("D24").Value = "'14.97"
$ws.Range("E24").Value = "  -1.15%  "
$ws.Range("E25").Value = "  +2.03%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").Value = "'9.28"
$ws.Range("E27").Value = "  +4.23%  "
$ws.Range("E28").Value = "  +3.16%  "
$ws.Range("E29").Value = "  +8.31%  "
$ws.Range("D30").Value = "'7.08"
$ws.Range("E30").Value = "  +13.84%  "
$ws.Range("D31").Value = "'2.94"
$ws.Range("E31").Value = "  +5.84%  "
$ws.Range("D32").Value = "'28.38"
$ws.Range("E32").Value = "  +2.94%  "
$ws.Range("E33").Value = "  +3.15%  "
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("E35").Value = "  +1.38%  "
$ws.Range("D36").Value = "'510.87"
$ws.Range("E36").Value = "  +5.46%  "
$ws.Range("D37").Value = "'54.88"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").Value = "'0.0904"
$ws.Range("E38").Value = "  +1.70%  "
$ws.Range("D39").Value = "'0.0426"
$ws.Range("E39").Value = "  +2.58%  "
$ws.Range("E40").Value = "  +10.26%  "
$ws.Range("E41").Value = "  +2.34%  "
$ws.Range("D42").Value = "'2.89"
$ws.Range("E42").Value = "  -0.68%  "
$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D43").Value = "0.0₃0679"
$ws.Range("E43").Value = "  +15.97%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.301"
$ws.Range("E44").Value = "  +7.12%  "
$ws.Range("E45").Value = "  +1.33%  "
$ws.Range("D46").Value = "2.901.45"
$ws.Range("E46").Value = "  -2.92%  "
$ws.Range("D47").Value = "'28.53"
$ws.Range("E47").Value = "  +1.79%  "
$ws.Range("E48").Value = "  +14.37%  "
$ws.Range("E49").Value = "  +3.83%  "
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("E51").Value = "  +5.71%  "
